$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('G2').Value = 'servinaz@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg'
$ws.Range('L6').Value = 60
$ws.Range('L7').Value = 22
$ws.Range('L9').Value = '34.1%'
$ws.Range('G10').Value = 'amany.raafat@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, esraa.mostafa@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg'
$ws.Range('L10').Value = '35.3%'
$ws.Range('G15').Value = 'hana.amr@med.asu.edu.eg, nancy.abdelshafy@med.asu.edu.eg'
$ws.Range('O15').Value = 9
$ws.Range('P15').Value = 2
$ws.Range('R15').Value = '40.9%'
$ws.Range('S15').Value = '37.2%'
$ws.Range('G18').Value = 'shorokmohamed@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg'
$ws.Range('G19').Value = 'wafaa.ebida@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg'
$ws.Range('G20').Value = 'Wafaa.ebida@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, marinasorial@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg'
$ws.Range('S22').Value = '20.6%'
$ws.Range('G24').Value = 'servinaz@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg'
$ws.Range('G25').Value = 'gehanadel@med.asu.edu.eg, alshimaa.atef@med.asu.edu.egm, manar.montaser@med.asu.edu.eg, backup@backdoor.com'
$ws.Range('G32').Value = 'amany.raafat@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, esraa.mostafa@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg'
$ws.Range('G37').Value = 'hana.amr@med.asu.edu.eg, nancy.abdelshafy@med.asu.edu.eg'
$ws.Range('G40').Value = 'shorokmohamed@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg'
$ws.Range('G41').Value = 'wafaa.ebida@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg'
$ws.Range('G42').Value = 'Wafaa.ebida@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, marinasorial@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg'
$ws.Range('G54').Value = 'yassmina.fattoh@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, maimustafa@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, merna.said@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg'
$ws.Range('G62').Value = 'shorokmohamed@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg'
$ws.Range('G63').Value = 'wafaa.ebida@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg'
$ws.Range('G76').Value = 'yassmina.fattoh@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, maimustafa@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, merna.said@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg'
$ws.Range('G81').Value = 'Walaa.h.ghanima@med.asu.edu.eg, enas.omran@med.asu.edu.eg, user@user.com'
$ws.Range('G84').Value = 'shorokmohamed@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg'
$ws.Range('G85').Value = 'wafaa.ebida@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg'
$ws.Range('G90').Value = 'shaimaa.ahmed@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg, manar.montaser@med.asu.edu.eg'
$ws.Range('G96').Value = 'mariam.noureldin@med.asu.edu.eg, norhan.mohamed@med.asu.edu.eg, Sara_nabil@med.asu.edu.eg, aml.awwad@med.asu.edu.eg'
$ws.Range('G98').Value = 'yassmina.fattoh@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, maimustafa@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, merna.said@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg'
$ws.Range('G106').Value = 'wafaa.ebida@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg'
$ws.Range('G107').Value = 'wafaa.ebida@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg'
$ws.Range('G108').Value = 'Wafaa.ebida@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, marinasorial@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg'
$ws.Range('G112').Value = 'shaimaa.ahmed@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg, manar.montaser@med.asu.edu.eg'
$ws.Range('G118').Value = 'mariam.noureldin@med.asu.edu.eg, norhan.mohamed@med.asu.edu.eg, Sara_nabil@med.asu.edu.eg, aml.awwad@med.asu.edu.eg'
$ws.Range('G120').Value = 'yassmina.fattoh@med.asu.edu.eg, amany.raafat@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, arwaelsayed03@med.asu.edu.eg, maimustafa@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, merna.said@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, maryam.ahmed@med.asu.edu.eg'
$ws.Range('G128').Value = 'wafaa.ebida@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg'
$ws.Range('G129').Value = 'wafaa.ebida@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg'
$ws.Range('G130').Value = 'Wafaa.ebida@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, marinasorial@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg'
$ws.Range('G134').Value = 'Amira.Sobhy@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg'
$ws.Range('G142').Value = 'yassmina.fattoh@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, esraa.mostafa@med.asu.edu.eg, merna.said@med.asu.edu.eg'
$ws.Range('G150').Value = 'wafaa.ebida@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg'
$ws.Range('G151').Value = 'wafaa.ebida@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg'
$ws.Range('G156').Value = 'manar.montaser@med.asu.edu.eg, alshimaa.atef@med.asu.edu.egm, Mohammedeltanany@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg'
$ws.Range('G164').Value = 'yassmina.fattoh@med.asu.edu.eg, marwa_mustafa@med.asu.edu.eg, basma.hamed@med.asu.edu.eg, esraa.mostafa@med.asu.edu.eg, merna.said@med.asu.edu.eg'
$ws.Range('G172').Value = 'wafaa.ebida@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg'
$ws.Range('G173').Value = 'wafaa.ebida@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg'
$ws.Range('H173').Value = '39/226'

# Row 21: status changes from "Not Recorded" to "Recorded" - copy formatting from a
# "Recorded" row (row 2) onto row 21, then set the updated values
$ws.Range("A2:I2").Copy()
$ws.Range("A21:I21").PasteSpecial(-4122)
$ws.Range("G21").Value = 'yasmin.m.senosy@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg'
$ws.Range("H21").Value = '37/216'
$ws.Range("I21").Value = 'Recorded'
